# Incorporo nuevos datos hasta diciembre de 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 82 already exists as a styled-but-empty placeholder row ---
# (A82 carries style s=2, D82 carries style s=3). Filling the values keeps
# those inherited styles for A82, but the real edit drops the date-style
# formatting that had been sitting on D82, so clear its format first.
$ws.Cells.Item(82, 1).Value = "041 Alquiler de vivienda"
$ws.Cells.Item(82, 2).Value = 2025
$ws.Cells.Item(82, 3).Value = 8
$ws.Cells.Item(82, 4).ClearFormats()
$ws.Cells.Item(82, 4).Value = 14.5
$ws.Cells.Item(82, 5).Value = 107.40740740740742

# --- New rows 83-86: copy column A's styling down (matches the "Componente"
# column formatting used throughout the table) before filling the values ---
$ws.Range("A82").Copy()
$ws.Range("A83:A86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(83, 1).Value = "041 Alquiler de vivienda"
$ws.Cells.Item(83, 2).Value = 2025
$ws.Cells.Item(83, 3).Value = 9
$ws.Cells.Item(83, 4).Value = 14.5
$ws.Cells.Item(83, 5).Value = 107.40740740740742

$ws.Cells.Item(84, 1).Value = "041 Alquiler de vivienda"
$ws.Cells.Item(84, 2).Value = 2025
$ws.Cells.Item(84, 3).Value = 10
$ws.Cells.Item(84, 4).Value = 14.5
$ws.Cells.Item(84, 5).Value = 107.40740740740742

$ws.Cells.Item(85, 1).Value = "041 Alquiler de vivienda"
$ws.Cells.Item(85, 2).Value = 2025
$ws.Cells.Item(85, 3).Value = 11
$ws.Cells.Item(85, 4).Value = 14.6
$ws.Cells.Item(85, 5).Value = 108.14814814814815

$ws.Cells.Item(86, 1).Value = "041 Alquiler de vivienda"
$ws.Cells.Item(86, 2).Value = 2025
$ws.Cells.Item(86, 3).Value = 12
$ws.Cells.Item(86, 4).Value = 14.7
$ws.Cells.Item(86, 5).Value = 108.88888888888889

# --- Match the author's final cursor/scroll position in the saved view ---
$ws.Range("F90").Select()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
